$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are plain (unformatted) text strings in the source file.
# Some new "Price" values are digit-strings that Excel would otherwise
# auto-convert to numbers (e.g. "0.7223"), so those cells are temporarily
# forced to Text format while assigning, then restored to the default style
# (matching the original, unstyled cells) so no numeric coercion occurs.

$ws.Range("D2").Value = "29.262.41"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.866.86"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7223"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07825"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3085"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08247"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7211"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "1.848.21"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "29.276.41"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.856"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007815"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "2.100.32"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.982"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.29%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.964"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.391"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05192"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.932"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7274"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01856"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.701"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "1.174.19"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9021"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.094"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5286"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "1.997.95"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.890"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.293"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
